$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) / Volume(1h) (column E) refresh, mirroring a data-source re-pull.
# Column D holds numeric-looking values as TEXT (e.g. "27.061.42", "0.0₃0740"), so
# any replacement that parses as a plain number is written with a leading "'" (quote
# prefix) to keep Excel from reinterpreting the cell as a Number.

$ws.Range("D2").Value = '27.053.40'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '1.673.97'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''215.07'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '''0.517'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +1.89%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").Value = '''21.17'
$ws.Range("E10").Value = '  +4.52%  '
$ws.Range("D11").Value = '''0.0882'
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("D12").Value = '1.911.15'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '1.678.88'
$ws.Range("E13").Value = '  +2.21%  '
$ws.Range("D14").Value = '''4.11'
$ws.Range("D15").Value = '''0.533'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '''65.96'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '27.042.13'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '''8.15'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("E19").Value = '  +1.52%  '
$ws.Range("D20").Value = '0.0₃0739'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '''4.45'
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("E23").Value = '  +2.03%  '
$ws.Range("E24").Value = '  -1.86%  '
$ws.Range("D25").Value = '''146.17'
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  +1.52%  '
$ws.Range("D27").Value = '''16.34'
$ws.Range("E27").Value = '  +2.59%  '
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").Value = '1.542.37'
$ws.Range("E33").Value = '  +5.91%  '
$ws.Range("E34").Value = '  +1.74%  '
$ws.Range("E35").Value = '  +3.54%  '
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").Value = '''0.917'
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("E40").Value = '  +2.56%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").Value = '''67.76'
$ws.Range("E42").Value = '  +2.23%  '
$ws.Range("D43").Value = '''5.59'
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("D44").Value = '''2.25'
$ws.Range("E44").Value = '  -1.98%  '
$ws.Range("D45").Value = '1.818.29'
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("D47").Value = '''90.84'
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").Value = '''8.03'
$ws.Range("E50").Value = '  +5.04%  '
$ws.Range("E51").Value = '  +0.57%  '
